$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "2021年" data row (row 5) below the existing 2018-2020 rows,
# growing the used range from A1:DK4 to A1:DK5.

# Copy row 4's label-cell style (bold, bordered, centered) onto A5, then
# overwrite the copied text with the new row's label.
$ws.Range("A4").Copy($ws.Range("A5"))

$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 48.1
$ws.Range("C5").Value = 69.90000000000001
$ws.Range("D5").Value = 34.4
$ws.Range("E5").Value = 172.8
$ws.Range("F5").Value = -7.9
$ws.Range("H5").Value = -29.3
$ws.Range("I5").Value = 45.6
$ws.Range("J5").Value = 34.3
$ws.Range("K5").Value = -45.4
$ws.Range("L5").Value = -34.7
$ws.Range("M5").Value = -16.3
$ws.Range("N5").Value = -94.5
$ws.Range("O5").Value = -21.5
$ws.Range("P5").Value = 9
$ws.Range("Q5").Value = -78.40000000000001
$ws.Range("R5").Value = 9.1
$ws.Range("T5").Value = -30.3
$ws.Range("U5").Value = -5.5
$ws.Range("V5").Value = -2.9
$ws.Range("W5").Value = -9.699999999999999
$ws.Range("X5").Value = 7.5
$ws.Range("Y5").Value = 34.3
$ws.Range("Z5").Value = 13.9
$ws.Range("AA5").Value = 131.1
$ws.Range("AB5").Value = 3.5
$ws.Range("AC5").Value = 3.7
$ws.Range("AD5").Value = -32.6
$ws.Range("AE5").Value = 3.2
$ws.Range("AF5").Value = -95.59999999999999
$ws.Range("AH5").Value = -68.8
$ws.Range("AI5").Value = -30.8
$ws.Range("AJ5").Value = 222.1
$ws.Range("AK5").Value = -14.6
$ws.Range("AL5").Value = 100.1
$ws.Range("AN5").Value = -66.7
$ws.Range("AO5").Value = -52.9
$ws.Range("AP5").Value = -21.5
$ws.Range("AQ5").Value = 2.6
$ws.Range("AR5").Value = -23.7
$ws.Range("AU5").Value = 36
$ws.Range("AV5").Value = 106.5
$ws.Range("AW5").Value = -47.1
$ws.Range("AX5").Value = -12.6
$ws.Range("AY5").Value = -12
$ws.Range("AZ5").Value = 45.9
$ws.Range("BA5").Value = 23.2
$ws.Range("BB5").Value = 46.8
$ws.Range("BC5").Value = -98
$ws.Range("BD5").Value = -41.8
$ws.Range("BE5").Value = -34.6
$ws.Range("BF5").Value = 77.8
$ws.Range("BG5").Value = -87
$ws.Range("BH5").Value = -99.59999999999999
$ws.Range("BI5").Value = -36.1
$ws.Range("BJ5").Value = -16.9
$ws.Range("BK5").Value = 9.699999999999999
$ws.Range("BL5").Value = -22.7
$ws.Range("BM5").Value = -40.8
$ws.Range("BN5").Value = -29.6
$ws.Range("BO5").Value = -27.2
$ws.Range("BP5").Value = -99.2
$ws.Range("BQ5").Value = 79.90000000000001
$ws.Range("BR5").Value = 16.5
$ws.Range("BS5").Value = 2.8
$ws.Range("BT5").Value = 6.8
$ws.Range("BU5").Value = -17.4
$ws.Range("BV5").Value = -7.9
$ws.Range("BW5").Value = -7.8
$ws.Range("BX5").Value = 4.2
$ws.Range("BY5").Value = -38.1
$ws.Range("BZ5").Value = -47.5
$ws.Range("CA5").Value = -33.5
$ws.Range("CB5").Value = -8.1
$ws.Range("CC5").Value = 40.1
$ws.Range("CD5").Value = -68.90000000000001
$ws.Range("CE5").Value = 5.6
$ws.Range("CF5").Value = 49.9
$ws.Range("CG5").Value = 120.6
$ws.Range("CH5").Value = 238.3
$ws.Range("CI5").Value = 59.3
$ws.Range("CJ5").Value = 34.5
$ws.Range("CK5").Value = 140.2
$ws.Range("CL5").Value = -73.90000000000001
$ws.Range("CM5").Value = 174.7
$ws.Range("CN5").Value = 26.7
$ws.Range("CO5").Value = 5
$ws.Range("CP5").Value = -5.4
$ws.Range("CQ5").Value = -3.5
$ws.Range("CR5").Value = 29.4
$ws.Range("CS5").Value = -16.7
$ws.Range("CT5").Value = 2.8
$ws.Range("CU5").Value = 37.4
$ws.Range("CV5").Value = -58.5
$ws.Range("CW5").Value = 182.2
$ws.Range("CX5").Value = -23.3
$ws.Range("CY5").Value = 4.4
$ws.Range("CZ5").Value = 45.9
$ws.Range("DA5").Value = -4.6
$ws.Range("DB5").Value = -8.9
$ws.Range("DC5").Value = -24.8
$ws.Range("DD5").Value = -14.8
$ws.Range("DE5").Value = -11.6
$ws.Range("DF5").Value = 10.2
$ws.Range("DG5").Value = 263.6
$ws.Range("DH5").Value = -30.9
$ws.Range("DI5").Value = -61.2
$ws.Range("DJ5").Value = 56.8
$ws.Range("DK5").Value = 9.800000000000001

# A handful of columns have no data point for 2021 (they are present in the
# row as empty text cells rather than simply missing). A leading apostrophe
# forces Excel to store an empty-text cell instead of clearing it back to
# blank; ClearFormats() then drops the transient "quote prefix" formatting
# so the cell is left with no explicit style, matching the other blank
# cells already in the sheet (e.g. G4, P4, R4).
$ws.Range("G5").Value = "'"
$ws.Range("G5").ClearFormats()
$ws.Range("S5").Value = "'"
$ws.Range("S5").ClearFormats()
$ws.Range("AG5").Value = "'"
$ws.Range("AG5").ClearFormats()
$ws.Range("AM5").Value = "'"
$ws.Range("AM5").ClearFormats()
$ws.Range("AS5").Value = "'"
$ws.Range("AS5").ClearFormats()
$ws.Range("AT5").Value = "'"
$ws.Range("AT5").ClearFormats()
